# CasosColombia.xlsx update
# - Swap a handful of stray "NaN" text markers for their correct numeric
#   values (and vice versa) in existing rows.
# - Append a new data row (186) for 2020-09-13 (serial 44080).
# - Move the frozen-pane selection to the new last cell, DX186.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections -----------------------------------------
$ws.Range("I14").Value  = 1
$ws.Range("CU17").Value = 1
$ws.Range("L18").Value  = "NaN"
$ws.Range("CK23").Value = 1
$ws.Range("CU27").Value = "NaN"
$ws.Range("AK34").Value = 1
$ws.Range("AY44").Value = "NaN"
$ws.Range("AY45").Value = 4
$ws.Range("O51").Value  = "NaN"
$ws.Range("O52").Value  = 14
$ws.Range("DN59").Value = 1
$ws.Range("DN61").Value = "NaN"
$ws.Range("AP83").Value = "NaN"
$ws.Range("W87").Value  = "NaN"
$ws.Range("AP87").Value = 5
$ws.Range("CK89").Value = "NaN"

# --- New row 186 (columns A:DX) -------------------------------------------
$row = 186
$rowValues = @(44080,666521,2714,88674,65074,225947,25967,4342,3396,6789,6397,12934,3786,20991,25829,5884,5991,13083,10031,15151,12744,3213,1572,7047,21705,12492,8349,50180,1274,305,448,453,171,135,351,1985,3706,36660,7556,2432,39336,999,20922,1473,8828,1559,1581,5369,1682,951,2481,2647,51319,12855,3428,8183,4692,280,1418,2626,734,2066,8766,8823,9142,13958,1892,835,9396,8011,9458,1789,1674,3820,3810,1132,4937,2725,1483,795,2364,2032,1488,1105,5497,1630,1247,1418,1803,1671,1991,1280,1113,1131,665,3117,1155,823,806,1435,1337,691,771,1040,1291,1129,1247,960,319,344,731,652,427,534,352,627,720,517,481,372,517,122942,282210,12288,122147,75514,34498,10173)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item($row, $i + 1).Value = $rowValues[$i]
}

# --- View: select the new bottom-right-most cell ---------------------------
$ws.Activate()
$ws.Range("DX186").Select()
